$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1423.25
$ws.Range("I43").Value = 3333
$ws.Range("J43").Value = 1150.4286
$ws.Range("K43").Value = 3333
$ws.Range("L43").Value = 1150.4286
$ws.Range("M43").Value = -3264
$ws.Range("N43").Value = -1288.4286
$ws.Range("H116").Value = 6441306
$ws.Range("I116").Value = 7085251.5
$ws.Range("J116").Value = 1850
$ws.Range("K116").Value = 7085251.5
$ws.Range("L116").Value = 1850
$ws.Range("M116").Value = -7081809.5
$ws.Range("N116").Value = -8734

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 19231568
$ws.Range("I2").Value = 23810224
$ws.Range("K2").Value = 23810224
$ws.Range("M2").Value = -23810111
$ws.Range("H28").Value = 3826.8333
$ws.Range("I28").Value = 992.2
$ws.Range("K28").Value = 992.2
$ws.Range("M28").Value = -800.2
$ws.Range("H31").Value = 5501.4443
$ws.Range("I31").Value = 5501.4443
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 5501.4443
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -5207.4443
$ws.Range("N31").Value = ""
$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("M44").Value = ""
$ws.Range("N44").Value = ""
$ws.Range("H61").Value = 1865.5
$ws.Range("I61").Value = 1619.0769
$ws.Range("J61").Value = 2933.3333
$ws.Range("K61").Value = 1619.0769
$ws.Range("L61").Value = 2933.3333
$ws.Range("M61").Value = -1407.0769
$ws.Range("N61").Value = -3357.3333
$ws.Range("H99").Value = 3826.8333
$ws.Range("I99").Value = 992.2
$ws.Range("K99").Value = 992.2
$ws.Range("M99").Value = 2002.8
$ws.Range("H102").Value = 90910900
$ws.Range("I102").Value = 142858270
$ws.Range("J102").Value = 3000
$ws.Range("K102").Value = 142858270
$ws.Range("L102").Value = 3000
$ws.Range("M102").Value = -142856648
$ws.Range("N102").Value = -6244
$ws.Range("H114").Value = 27560
$ws.Range("J114").Value = 27560
$ws.Range("L114").Value = 27560
$ws.Range("N114").Value = -36238
$ws.Range("H116").Value = 19231568
$ws.Range("I116").Value = 23810224
$ws.Range("K116").Value = 23810224
$ws.Range("M116").Value = -23807930
$ws.Range("H136").Value = 1865.5
$ws.Range("I136").Value = 1619.0769
$ws.Range("J136").Value = 2933.3333
$ws.Range("K136").Value = 4857.2307
$ws.Range("L136").Value = 8799.999899999999
$ws.Range("M136").Value = -2307.2307
$ws.Range("N136").Value = -13899.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 19231568
$ws.Range("I3").Value = 23810224
$ws.Range("K3").Value = 23810224
$ws.Range("M3").Value = -23810110
$ws.Range("H94").Value = 8678.4
$ws.Range("I94").Value = 708.4091
$ws.Range("J94").Value = 67125
$ws.Range("K94").Value = 708.4091
$ws.Range("L94").Value = 67125
$ws.Range("M94").Value = -257.4091
$ws.Range("N94").Value = -68027
$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("M102").Value = ""
$ws.Range("H105").Value = 2152.353
$ws.Range("I105").Value = 1761.7273
$ws.Range("J105").Value = 2868.5
$ws.Range("K105").Value = 1761.7273
$ws.Range("L105").Value = 2868.5
$ws.Range("M105").Value = -14.72730000000001
$ws.Range("N105").Value = -6362.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 505.46667
$ws.Range("I22").Value = 530
$ws.Range("J22").Value = 493.2
$ws.Range("K22").Value = 530
$ws.Range("L22").Value = 493.2
$ws.Range("M22").Value = -180
$ws.Range("N22").Value = -1193.2
$ws.Range("H31").Value = 6575038
$ws.Range("I31").Value = 4675361
$ws.Range("J31").Value = 11113155
$ws.Range("K31").Value = 4675361
$ws.Range("L31").Value = 11113155
$ws.Range("M31").Value = -4675066
$ws.Range("N31").Value = -11113745
$ws.Range("H34").Value = 6575038
$ws.Range("I34").Value = 4675361
$ws.Range("J34").Value = 11113155
$ws.Range("K34").Value = 4675361
$ws.Range("L34").Value = 11113155
$ws.Range("M34").Value = -4675159
$ws.Range("N34").Value = -11113559

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 42110464
$ws.Range("I70").Value = 66671584
$ws.Range("J70").Value = 5685.7144
$ws.Range("K70").Value = 66671584
$ws.Range("L70").Value = 5685.7144
$ws.Range("M70").Value = -66671314
$ws.Range("N70").Value = -6225.7144
$ws.Range("H73").Value = 42110464
$ws.Range("I73").Value = 66671584
$ws.Range("J73").Value = 5685.7144
$ws.Range("K73").Value = 66671584
$ws.Range("L73").Value = 5685.7144
$ws.Range("M73").Value = -66670648
$ws.Range("N73").Value = -7557.7144
$ws.Range("H99").Value = 4078
$ws.Range("I99").Value = 1467.75
$ws.Range("J99").Value = 14519
$ws.Range("K99").Value = 1467.75
$ws.Range("L99").Value = 14519
$ws.Range("M99").Value = 778.25
$ws.Range("N99").Value = -19011
$ws.Range("H102").Value = 1817.6538
$ws.Range("I102").Value = 1852.95
$ws.Range("K102").Value = 1852.95
$ws.Range("M102").Value = -230.95
$ws.Range("H113").Value = 1055.3549
$ws.Range("I113").Value = 840.82355
$ws.Range("J113").Value = 1315.8572
$ws.Range("K113").Value = 840.82355
$ws.Range("L113").Value = 1315.8572
$ws.Range("M113").Value = 1329.17645
$ws.Range("N113").Value = -5655.8572
$ws.Range("H132").Value = 2520.0625
$ws.Range("I132").Value = 1847.6364
$ws.Range("J132").Value = 3999.4
$ws.Range("K132").Value = 5542.9092
$ws.Range("L132").Value = 11998.2
$ws.Range("M132").Value = -3012.9092
$ws.Range("N132").Value = -17058.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 277.27777
$ws.Range("I22").Value = 249
$ws.Range("J22").Value = 299.9
$ws.Range("K22").Value = 249
$ws.Range("L22").Value = 299.9
$ws.Range("M22").Value = 46
$ws.Range("N22").Value = -889.9
$ws.Range("H27").Value = 277.27777
$ws.Range("I27").Value = 249
$ws.Range("J27").Value = 299.9
$ws.Range("K27").Value = 249
$ws.Range("L27").Value = 299.9
$ws.Range("M27").Value = -142
$ws.Range("N27").Value = -513.9
$ws.Range("H46").Value = 926.4583
$ws.Range("I46").Value = 499.75
$ws.Range("J46").Value = 1353.1666
$ws.Range("K46").Value = 499.75
$ws.Range("L46").Value = 1353.1666
$ws.Range("M46").Value = -311.75
$ws.Range("N46").Value = -1729.1666
$ws.Range("H55").Value = 147.5
$ws.Range("H88").Value = 29000
$ws.Range("J88").Value = 29000
$ws.Range("L88").Value = 29000
$ws.Range("N88").Value = -29856
$ws.Range("H91").Value = 29000
$ws.Range("J91").Value = 29000
$ws.Range("L91").Value = 29000
$ws.Range("N91").Value = -31964
$ws.Range("H100").Value = 1545.8182
$ws.Range("I100").Value = 1228.5714
$ws.Range("J100").Value = 2101
$ws.Range("K100").Value = 1228.5714
$ws.Range("L100").Value = 2101
$ws.Range("M100").Value = -687.5714
$ws.Range("N100").Value = -3183

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").Value = ""
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").Value = ""
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("M69").Value = ""
$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("M72").Value = ""
$ws.Range("H113").Value = 27778870
$ws.Range("I113").Value = 33334316
$ws.Range("J113").Value = 1633.3334
$ws.Range("K113").Value = 100002948
$ws.Range("L113").Value = 4900.0002
$ws.Range("M113").Value = -100000778
$ws.Range("N113").Value = -9240.0002
